$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Cells.Item(1,2).Value = "FIPS"
$ws.Cells.Item(1,3).Value = "Year"
$ws.Cells.Item(1,4).Value = "Population"

# New header cells E1 / F1 need the same style as the existing header cells.
# Copy formatting from B1 (which already carries the bold/border/center style)
# onto E1:F1 before writing their values so no new style entry is created.
$ws.Range("B1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Cells.Item(1,5).Value = "Deaths Per 100,000 People"
$ws.Cells.Item(1,6).Value = "CalculatedDeaths"

# --- Row 2 (count) ---
$ws.Cells.Item(2,2).Value = 39819
$ws.Cells.Item(2,3).Value = 39819
$ws.Cells.Item(2,4).Value = 39819
$ws.Cells.Item(2,5).Value = 39819
$ws.Cells.Item(2,6).Value = 39819

# --- Row 3 (mean) ---
$ws.Cells.Item(3,2).Value = 30511.0365654587
$ws.Cells.Item(3,3).Value = 2009
$ws.Cells.Item(3,4).Value = 100070.0092920465
$ws.Cells.Item(3,5).Value = 12.5066937806446
$ws.Cells.Item(3,6).Value = 12.38052186142294

# --- Row 4 (std) ---
$ws.Cells.Item(4,2).Value = 15086.75968991359
$ws.Cells.Item(4,3).Value = 3.741704370975745
$ws.Cells.Item(4,4).Value = 316386.0302388998
$ws.Cells.Item(4,5).Value = 8.228459764040865
$ws.Cells.Item(4,6).Value = 38.16859162489131

# --- Row 5 (min) ---
$ws.Cells.Item(5,2).Value = 1001
$ws.Cells.Item(5,3).Value = 2003
$ws.Cells.Item(5,4).Value = 55
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 0

# --- Row 6 (25%) ---
$ws.Cells.Item(6,2).Value = 19007
$ws.Cells.Item(6,3).Value = 2006
$ws.Cells.Item(6,4).Value = 11734.5
$ws.Cells.Item(6,5).Value = 8.227687541724906
$ws.Cells.Item(6,6).Value = 1

# --- Row 7 (50%) ---
$ws.Cells.Item(7,2).Value = 29151
$ws.Cells.Item(7,3).Value = 2009
$ws.Cells.Item(7,4).Value = 26474
$ws.Cells.Item(7,5).Value = 11.50196491900699
$ws.Cells.Item(7,6).Value = 3

# --- Row 8 (75%) ---
$ws.Cells.Item(8,2).Value = 46033
$ws.Cells.Item(8,3).Value = 2012
$ws.Cells.Item(8,4).Value = 68041
$ws.Cells.Item(8,5).Value = 15.71633182457846
$ws.Cells.Item(8,6).Value = 9

# --- Row 9 (max) ---
$ws.Cells.Item(9,2).Value = 56045
$ws.Cells.Item(9,3).Value = 2015
$ws.Cells.Item(9,4).Value = 10085416
$ws.Cells.Item(9,5).Value = 126.5522423475441
$ws.Cells.Item(9,6).Value = 862
